$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "nickname" -> "charge"
$ws.Range("C1").Value = "charge"

# Update values in column C from nicknames to charge amounts (stored as text)
$ws.Range("C2").Value = "80"
$ws.Range("C3").Value = "40"
